$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Daily practice")
$ws3 = $wb.Worksheets.Item("Learnings")

# --- 1. Row 20 (Unit array): expand the mistake note with the post-editorial reflection ---
$ws1.Range("E20").Value = 'I was trying very hard to find out a condition where if number of negative is greater than or equal to number of positives, then how many negatives to be converted to get the answer. But I was always getting stuck with one test case or another.
Finally after hours of finding out cases, I got the solution myself (Will be checking editorial for sure)
Cases:
np --> number of 1s
nn --> number of -1s 
1. If np == n: return 0 (since all are positive)
2. if np>nn: return 0 if even number of -1s else return 1 (since 1s > -1s, sum will always be greater than 0, we need to take care of the product and if there are even number of -1s, product will also be =ve, so returning 0, else make 1 -1 as 1, so one operation)
3. This is the difficult part where nn>=np (see code in the column aside)
Here is the change of perspective from the editorial.
Now, since this is not a YES or NO question, they went ahead and performed the operations in a loop. That makes the code much easier. Yes, my code is O(1) to find the solution but to get the number of 1s and -1s I am anyways doing a O(n) so that doesnt help !
So, based on the solution given in the editorial, this indeed is a easy problem'

# --- 2. Learnings sheet: first new takeaway (tied to the Unit array problem) ---
$ws3.Range("A3").Value = 'Sometimes, when the question asks for a direct answer and the solution looks too mathematical, that is, a lot of equations and if-else complexities, try solving the question by doing the action itself since the action itself is very trivial. Complexity optimization takes the back bench then'

# --- 3. Row 21: Twin permutations (20th May, 2025) ---
$ws1.Range("A17:E17").Copy()
$ws1.Range("A21:E21").PasteSpecial(-4122)
$ws1.Range("A21").Value = '20th May, 2025'
$ws1.Range("B21").Value = 'Twin permutations'
$ws1.Range("C21").Value = 'https://codeforces.com/problemset/problem/1831/A'
$ws1.Range("D21").Value = "Difficult"
$ws1.Range("E21").Value = 'I am thinking of how to find that permutation.
Thought of sorting the array, failed for  1,2,4,5,3
Though of sorting the part which was unsorted. Failed for 1,3,8,5,6,10
Thinking of an idea, where I find the sum of max+min and then develop the array by subtracting the sum with every element. This works !'
$ws1.Rows.Item(21).RowHeight = 72

# --- 4. Row 22: Blank space (23rd May, 2025) ---
$ws1.Range("A19:E19").Copy()
$ws1.Range("A22:E22").PasteSpecial(-4122)
$ws1.Range("A22").Value = '23rd May, 2025'
$ws1.Range("B22").Value = 'Blank space'
$ws1.Range("C22").Value = 'https://codeforces.com/problemset/problem/1829/B'
$ws1.Range("D22").Value = "Easy"
$ws1.Range("E22").Value = 'Just keep track of the last one'
$ws1.Range("E2").Copy()
$ws1.Range("A22").PasteSpecial(-4122)

# --- 5. Row 23: Coins (23rd May, 2025) ---
$ws1.Range("A19:E19").Copy()
$ws1.Range("A23:E23").PasteSpecial(-4122)
$ws1.Range("A23").Value = '23rd May, 2025'
$ws1.Range("B23").Value = 'Coins'
$ws1.Range("C23").Value = 'https://codeforces.com/problemset/problem/1814/A'
$ws1.Range("D23").Value = "Easy"
$ws1.Range("E23").Value = 'Easy only because I got the observation correct
So, I had to find a way where the condition will always be false and after playing around with a few cases, it seemed that the condition will always be false if k is even but n is false
I have later gone through the proof of why this is always true. See notes'
$ws1.Range("E2").Copy()
$ws1.Range("A23").PasteSpecial(-4122)
$ws1.Rows.Item(23).RowHeight = 72

# --- 6. Row 24: Walking master (23rd May, 2025) ---
$ws1.Range("A18:E18").Copy()
$ws1.Range("A24:E24").PasteSpecial(-4122)
$ws1.Range("A24").Value = '23rd May, 2025'
$ws1.Range("B24").Value = 'Walking master'
$ws1.Range("C24").Value = 'https://codeforces.com/problemset/problem/1806/A'
$ws1.Range("D24").Value = "Difficult"
$ws1.Range("E24").Value = 'I found out 2 things:
1. If the point (c,d) is below the diagonal of (a,b) its not reachable. And we can find by the line equations
2. If the line is on the diagonal then we are good by simply subtracting (x1-x2)
But I am stuck that if its above the diagonal, how to determine ?
I found out that you can simply find out the exact point below it and to go right above its 2*(direct distance)
But this is not working
I then thought hard and worked out some maths
We need to move from (x1,y1) to (x2,y2)
If y2 < y1, we can never go with the left and diag moves
Also, if (x2,y2) is below the diagonal. i.e, y2<y1+x2-x1 then not possible.
Now, comes the case where its possible.
To do that we need to find a point (p,q) so that we can do k number of left moves and z number of diag moves.
Now, if you observe carefully, for (p,q) q will always be achieved by moving left, thus q = y1
Now p will be in the same diagonal as (x2,y2) 
Thus, p = y1 - y2 + x2
So, we get (p,q)
Thus, the total distance will be abs(x1 - p) + abs(y2 - q)'
$ws1.Range("E2").Copy()
$ws1.Range("A24").PasteSpecial(-4122)
$ws1.Rows.Item(24).RowHeight = 345.6

# --- 7. Learnings sheet: second new takeaway (tied to the Walking master problem) ---
$ws3.Range("A4").Value = 'Maths is essential for CF problems'

# --- 8. Column width on Learnings sheet A (wide single note column) ---
$ws3.Columns.Item(1).ColumnWidth = 254.9

# --- 9. Selections / active sheet, matching the saved view state ---
$ws1.Range("D24").Select()
$ws3.Range("A5").Select()
